$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 11.479533064867338
$ws.Range("C2").Value = 11.344356891707601
$ws.Range("D2").Value = 12.44210363328866
$ws.Range("E2").Value = 11.503110239949237

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 10.896010125357197
$ws.Range("C3").Value = 9.4093520939954178
$ws.Range("D3").Value = 11.044606909410176
$ws.Range("E3").Value = 11.266127978668678

# Update the selection to match the new sqref B1:E3
$ws.Range("B1:E3").Select()
